# "removed created by in import jurnal"
#
# The journal-import template had an extra "created_by" column (H) that is
# no longer used. This script removes that column, renames two of the
# remaining rows' descriptions, and updates the debit/credit amounts for
# the sample journal entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-unused "created_by" column (H1 header + H2/H3 data cells),
# shrinking the used range from A1:H3 down to A1:G3.
$ws.Range("H1:H3").Delete()

# Row 2 (debit line): update the description and the debit amount.
$ws.Range("D2").Value = "Pembelian AC"
$ws.Range("F2").Value = 1500000

# Row 3 (credit line): update the description and the credit amount.
$ws.Range("D3").Value = "Bank BNI"
$ws.Range("G3").Value = 1500000

# Match the saved view's active selection.
$null = $ws.Range("G9").Select()
